$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C, rows 2 through 27 hold a "Förändrad" (changed) date stored as the
# serial number 45279 (2023-12-19). Bump each of these to 45280 (2023-12-20),
# matching the rest of the row/cell (style, type, etc.) unchanged.
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -eq 45279) {
        $cell.Value = 45280
    }
}
